$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 168.5
$ws.Range("I8").Value = 152.2
$ws.Range("K8").Value = 456.6
$ws.Range("M8").Value = -317.6
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H137").Value = 3291.8704
$ws.Range("J137").Value = 4986.278
$ws.Range("L137").Value = 14958.834
$ws.Range("N137").Value = -20058.834
$ws.Range("H138").Value = 7028.964
$ws.Range("I138").Value = 1661.25
$ws.Range("K138").Value = 4983.75
$ws.Range("M138").Value = 156.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1693.0322
$ws.Range("I2").Value = 1325.56
$ws.Range("J2").Value = 3224.1667
$ws.Range("K2").Value = 1325.56
$ws.Range("L2").Value = 3224.1667
$ws.Range("M2").Value = -1212.56
$ws.Range("N2").Value = -3450.1667
$ws.Range("H32").Value = 27030118
$ws.Range("I32").Value = 29413436
$ws.Range("K32").Value = 29413436
$ws.Range("M32").Value = -29413149
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").ClearContents()
$ws.Range("H102").Value = 202605.9
$ws.Range("I102").Value = 287294.28
$ws.Range("K102").Value = 287294.28
$ws.Range("M102").Value = -285672.28
$ws.Range("H116").Value = 1693.0322
$ws.Range("I116").Value = 1325.56
$ws.Range("J116").Value = 3224.1667
$ws.Range("K116").Value = 1325.56
$ws.Range("L116").Value = 3224.1667
$ws.Range("M116").Value = 968.4400000000001
$ws.Range("N116").Value = -7812.1667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1693.0322
$ws.Range("I3").Value = 1325.56
$ws.Range("J3").Value = 3224.1667
$ws.Range("K3").Value = 1325.56
$ws.Range("L3").Value = 3224.1667
$ws.Range("M3").Value = -1211.56
$ws.Range("N3").Value = -3452.1667
$ws.Range("H94").Value = 785.9211
$ws.Range("I94").Value = 514.28125
$ws.Range("K94").Value = 514.28125
$ws.Range("M94").Value = -63.28125
$ws.Range("H105").Value = 6984.552
$ws.Range("I105").Value = 8518
$ws.Range("K105").Value = 8518
$ws.Range("M105").Value = -6771

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 18573.74
$ws.Range("I2").Value = 47.88889
$ws.Range("K2").Value = 287.33334
$ws.Range("M2").Value = -174.33334
$ws.Range("H68").Value = 3946.7693
$ws.Range("I68").Value = 4534.5
$ws.Range("K68").Value = 13603.5
$ws.Range("M68").Value = -12792.5
$ws.Range("H71").Value = 3946.7693
$ws.Range("I71").Value = 4534.5
$ws.Range("K71").Value = 40810.5
$ws.Range("M71").Value = -36754.5
$ws.Range("H131").Value = 11537703
$ws.Range("I131").Value = 80826
$ws.Range("J131").Value = 20846416
$ws.Range("K131").Value = 242478
$ws.Range("L131").Value = 62539248
$ws.Range("M131").Value = -237438
$ws.Range("N131").Value = -62549328

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5142.2
$ws.Range("I22").Value = 5237
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 5237
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -4942
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 5142.2
$ws.Range("I27").Value = 5237
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 5237
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -5130
$ws.Range("N27").Value = -5214
$ws.Range("H61").Value = 4809.625
$ws.Range("I61").Value = 4083.4443
$ws.Range("J61").Value = 5743.2856
$ws.Range("K61").Value = 4083.4443
$ws.Range("L61").Value = 5743.2856
$ws.Range("M61").Value = -3881.4443
$ws.Range("N61").Value = -6147.2856
$ws.Range("H63").Value = 69590.60000000001
$ws.Range("J63").Value = 62999
$ws.Range("L63").Value = 62999
$ws.Range("N63").Value = -64497
$ws.Range("H66").Value = 69590.60000000001
$ws.Range("J66").Value = 62999
$ws.Range("L66").Value = 188997
$ws.Range("N66").Value = -196485
$ws.Range("H113").Value = 4809.625
$ws.Range("I113").Value = 4083.4443
$ws.Range("J113").Value = 5743.2856
$ws.Range("K113").Value = 4083.4443
$ws.Range("L113").Value = 5743.2856
$ws.Range("M113").Value = -1913.4443
$ws.Range("N113").Value = -10083.2856
$ws.Range("H136").Value = 1179434
$ws.Range("I136").Value = 1431276.6
$ws.Range("K136").Value = 4293829.800000001
$ws.Range("M136").Value = -4291279.800000001
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 40400736
$ws.Range("J100").Value = 947.4286
$ws.Range("L100").Value = 1894.8572
$ws.Range("N100").Value = -2976.8572
$ws.Range("H107").Value = 1558.4584
$ws.Range("I107").Value = 570.6
$ws.Range("J107").Value = 3204.889
$ws.Range("K107").Value = 1711.8
$ws.Range("L107").Value = 9614.667000000001
$ws.Range("M107").Value = 208.1999999999998
$ws.Range("N107").Value = -13454.667
$ws.Range("H113").Value = 647.3103599999999
$ws.Range("I113").Value = 556.2174
$ws.Range("J113").Value = 996.5
$ws.Range("K113").Value = 1668.6522
$ws.Range("L113").Value = 2989.5
$ws.Range("M113").Value = 501.3478
$ws.Range("N113").Value = -7329.5
$ws.Range("H123").Value = 109054.14
$ws.Range("J123").Value = 109054.14
$ws.Range("L123").Value = 109054.14
$ws.Range("N123").Value = -118854.14
$ws.Range("H125").Value = 90000
$ws.Range("J125").Value = 90000
$ws.Range("L125").Value = 90000
$ws.Range("N125").Value = -99840
$ws.Range("H132").Value = 1681.2858
$ws.Range("I132").Value = 1620.8064
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 4862.4192
$ws.Range("L132").Value = 6450
$ws.Range("M132").Value = -2332.4192
$ws.Range("N132").Value = -11510
